# doc(general): change releases location
#
# The paragraph that used to read:
#   "Созданные архивы размещаются в <bookmark/>https://github.com/oramake/oramake-framework/releases"
# is split into two paragraphs right before the (zero-width) "_GoBack"
# bookmark that precedes the URL, so the link (and the bookmark) move
# into their own paragraph. A trailing "." and a trailing space are
# appended after the URL, and the paragraph that used to follow it (the
# "Замечание: ... sourceforge ... специальные права." remark about
# sourceforge permissions) is removed entirely.

$d = $word.ActiveDocument

# Locate the zero-width "_GoBack" bookmark that sits right before the
# "https://github.com/..." URL run — it marks exactly where the new
# paragraph break has to go.
$bm = $d.Bookmarks.Item("_GoBack")

# Split the paragraph right before the bookmark/URL: everything up to
# and including "в " stays in the first paragraph; the bookmark + URL
# move into a brand-new paragraph that inherits the same paragraph
# formatting (ind/jc) via the split.
$splitRange = $d.Range($bm.Start, $bm.Start)
$splitRange.InsertParagraphAfter()

# Re-resolve the bookmark after the split and find the paragraph that
# now starts with it (i.e. the paragraph now holding the URL).
$bm = $d.Bookmarks.Item("_GoBack")
$count = $d.Paragraphs.Count
$urlParagraph = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -eq $bm.Start) {
        $urlParagraph = $p
        break
    }
}

# The very next paragraph is the old "Замечание: ... sourceforge ...
# специальные права." remark — remove it completely (text + paragraph
# mark), since the note no longer applies once the link gets its own
# paragraph.
$remarkParagraph = $urlParagraph.Next()
$remarkParagraph.Range.Delete()

# Append ". " (as two separate runs, matching the document's normal
# run-splitting behaviour) right after the URL, before the paragraph
# mark of its (now final) paragraph.
$urlParagraph.Range.InsertAfter(".")
$urlParagraph.Range.InsertAfter(" ")
